# Sprint 43 test case report - Day 2 summary counts entered
# and cursor moved to C11 after data entry (matches authored diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Day 2" test case summary numbers (rows 9-11).
$ws.Range("C9").Value = 6980
$ws.Range("C10").Value = 2278
$ws.Range("C11").Value = 2278

# Re-establish the header merges in row order so the saved file lists
# them sorted the same way Excel normalizes mergeCells on save.
$ws.Range("B2:C2").UnMerge()
$ws.Range("B8:C8").UnMerge()
$ws.Range("B14:C14").UnMerge()
$ws.Range("B20:C20").UnMerge()
$ws.Range("B26:C26").UnMerge()
$ws.Range("B32:C32").UnMerge()
$ws.Range("B38:C38").UnMerge()
$ws.Range("B44:C44").UnMerge()
$ws.Range("B50:C50").UnMerge()
$ws.Range("B57:C57").UnMerge()

$ws.Range("B2:C2").Merge()
$ws.Range("B8:C8").Merge()
$ws.Range("B14:C14").Merge()
$ws.Range("B20:C20").Merge()
$ws.Range("B26:C26").Merge()
$ws.Range("B32:C32").Merge()
$ws.Range("B38:C38").Merge()
$ws.Range("B44:C44").Merge()
$ws.Range("B50:C50").Merge()
$ws.Range("B57:C57").Merge()

# Leave the active selection on C11, matching where data entry ended.
$ws.Range("C11").Select()
